$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
